$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 611.25
$ws.Range("I9").Value = 611.25
$ws.Range("K9").Value = 611.25
$ws.Range("M9").Value = -442.25
$ws.Range("H86").Value = 2635527.8
$ws.Range("I86").Value = 3375.3333
$ws.Range("K86").Value = 3375.3333
$ws.Range("M86").Value = -2252.3333
$ws.Range("H89").Value = 2635527.8
$ws.Range("I89").Value = 3375.3333
$ws.Range("K89").Value = 16876.6665
$ws.Range("M89").Value = -11260.6665
$ws.Range("H125").Value = 6972.875
$ws.Range("J125").Value = 8333.333000000001
$ws.Range("L125").Value = 74999.997
$ws.Range("N125").Value = -79919.997
$ws.Range("H132").Value = 4401.5137
$ws.Range("I132").Value = 4384.722
$ws.Range("K132").Value = 13154.166
$ws.Range("M132").Value = -10624.166
$ws.Range("H138").Value = 4451.727
$ws.Range("J138").Value = 5451.933
$ws.Range("L138").Value = 16355.799
$ws.Range("N138").Value = -26635.799

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 996.7105
$ws.Range("I32").Value = 996.61646
$ws.Range("K32").Value = 996.61646
$ws.Range("M32").Value = -709.61646
$ws.Range("H45").Value = 85674.336
$ws.Range("I45").Value = 102009.2
$ws.Range("K45").Value = 102009.2
$ws.Range("M45").Value = -101632.2
$ws.Range("H97").Value = 1411.9375
$ws.Range("I97").Value = 1385.1072
$ws.Range("J97").Value = 1599.75
$ws.Range("K97").Value = 1385.1072
$ws.Range("L97").Value = 1599.75
$ws.Range("M97").Value = -889.1071999999999
$ws.Range("N97").Value = -2591.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1071.4445
$ws.Range("I16").Value = 994.2143
$ws.Range("K16").Value = 994.2143
$ws.Range("M16").Value = -707.2143
$ws.Range("H31").Value = 726282.9399999999
$ws.Range("I31").Value = 1304860.5
$ws.Range("J31").Value = 42509.547
$ws.Range("K31").Value = 1304860.5
$ws.Range("L31").Value = 42509.547
$ws.Range("M31").Value = -1304565.5
$ws.Range("N31").Value = -43099.547
$ws.Range("H34").Value = 726282.9399999999
$ws.Range("I34").Value = 1304860.5
$ws.Range("J34").Value = 42509.547
$ws.Range("K34").Value = 1304860.5
$ws.Range("L34").Value = 42509.547
$ws.Range("M34").Value = -1304658.5
$ws.Range("N34").Value = -42913.547
$ws.Range("H50").Value = 23460.445
$ws.Range("J50").Value = 25018
$ws.Range("L50").Value = 25018
$ws.Range("N50").Value = -26268
$ws.Range("H51").Value = 14898.2
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 35775
$ws.Range("H60").Value = 19999
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 14898.2
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 79921
$ws.Range("I68").Value = 79942
$ws.Range("J68").Value = 79900
$ws.Range("K68").Value = 79942
$ws.Range("L68").Value = 79900
$ws.Range("M68").Value = -79193
$ws.Range("N68").Value = -81398
$ws.Range("H71").Value = 79921
$ws.Range("I71").Value = 79942
$ws.Range("J71").Value = 79900
$ws.Range("K71").Value = 239826
$ws.Range("L71").Value = 239700
$ws.Range("M71").Value = -236082
$ws.Range("N71").Value = -247188
$ws.Range("H107").Value = 1348.7646
$ws.Range("I107").Value = 969.75
$ws.Range("J107").Value = 2258.4
$ws.Range("K107").Value = 969.75
$ws.Range("L107").Value = 2258.4
$ws.Range("M107").Value = 950.25
$ws.Range("N107").Value = -6098.4
$ws.Range("H113").Value = 1071.4445
$ws.Range("I113").Value = 994.2143
$ws.Range("K113").Value = 994.2143
$ws.Range("M113").Value = 1175.7857
$ws.Range("H132").Value = 3541.0625
$ws.Range("I132").Value = 2204.3845
$ws.Range("J132").Value = 9333.333000000001
$ws.Range("K132").Value = 6613.1535
$ws.Range("L132").Value = 27999.999
$ws.Range("M132").Value = -4083.1535
$ws.Range("N132").Value = -33059.999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5919431.5
$ws.Range("J113").Value = 43301.2
$ws.Range("L113").Value = 129903.6
$ws.Range("N113").Value = -134243.6
$ws.Range("H136").Value = 5811.6
$ws.Range("J136").Value = 7332.6665
$ws.Range("L136").Value = 21997.9995
$ws.Range("N136").Value = -32197.9995
$ws.Range("H138").Value = 6549.6113
$ws.Range("I138").Value = 3177.25
$ws.Range("J138").Value = 9247.5
$ws.Range("K138").Value = 9531.75
$ws.Range("L138").Value = 27742.5
$ws.Range("M138").Value = -4391.75
$ws.Range("N138").Value = -38022.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.5
$ws.Range("I2").Value = 194.33333
$ws.Range("K2").Value = 194.33333
$ws.Range("M2").Value = -81.33332999999999
$ws.Range("H80").Value = 2005050
$ws.Range("I80").Value = 1430213.4
$ws.Range("J80").Value = 3346335.2
$ws.Range("K80").Value = 1430213.4
$ws.Range("L80").Value = 3346335.2
$ws.Range("M80").Value = -1429215.4
$ws.Range("N80").Value = -3348331.2
$ws.Range("H83").Value = 2005050
$ws.Range("I83").Value = 1430213.4
$ws.Range("J83").Value = 3346335.2
$ws.Range("K83").Value = 7151067
$ws.Range("L83").Value = 16731676
$ws.Range("M83").Value = -7146075
$ws.Range("N83").Value = -16741660
$ws.Range("H97").Value = 594.15625
$ws.Range("I97").Value = 651.1667
$ws.Range("J97").Value = 423.125
$ws.Range("K97").Value = 651.1667
$ws.Range("L97").Value = 423.125
$ws.Range("M97").Value = -155.1667
$ws.Range("N97").Value = -1415.125
$ws.Range("H123").Value = 38570.855
$ws.Range("J123").Value = 38570.855
$ws.Range("L123").Value = 38570.855
$ws.Range("N123").Value = -43470.855

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 426656.72
$ws.Range("I7").Value = 671850.5600000001
$ws.Range("K7").Value = 671850.5600000001
$ws.Range("M7").Value = -671738.5600000001
$ws.Range("H46").Value = 2508.7104
$ws.Range("I46").Value = 2396.7
$ws.Range("J46").Value = 2633.1667
$ws.Range("K46").Value = 2396.7
$ws.Range("L46").Value = 2633.1667
$ws.Range("M46").Value = -2208.7
$ws.Range("N46").Value = -3009.1667
$ws.Range("H68").Value = 62548.293
$ws.Range("J68").Value = 94639.09
$ws.Range("L68").Value = 94639.09
$ws.Range("N68").Value = -96137.09
$ws.Range("H71").Value = 62548.293
$ws.Range("J71").Value = 94639.09
$ws.Range("L71").Value = 473195.45
$ws.Range("N71").Value = -480683.45
$ws.Range("H93").Value = 1564.2632
$ws.Range("I93").Value = 1633.4445
$ws.Range("J93").Value = 319
$ws.Range("K93").Value = 1633.4445
$ws.Range("L93").Value = 319
$ws.Range("M93").Value = -385.4445000000001
$ws.Range("N93").Value = -2815
$ws.Range("H100").Value = 63849.668
$ws.Range("J100").Value = 17875
$ws.Range("L100").Value = 17875
$ws.Range("N100").Value = -18957
$ws.Range("H126").Value = 426656.72
$ws.Range("I126").Value = 671850.5600000001
$ws.Range("K126").Value = 2015551.68
$ws.Range("M126").Value = -2013081.68
$ws.Range("H136").Value = 8265813
$ws.Range("I136").Value = 16020116
$ws.Range("K136").Value = 48060348
$ws.Range("M136").Value = -48057798

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 70345.92999999999
$ws.Range("I96").Value = 127124.625
$ws.Range("K96").Value = 127124.625
$ws.Range("M96").Value = -125751.625
$ws.Range("H107").Value = 67553.87
$ws.Range("I107").Value = 111923.11
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 335769.33
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -333849.33
$ws.Range("N107").Value = -6840
$ws.Range("H136").Value = 277569.28
$ws.Range("J136").Value = 304572.16
$ws.Range("L136").Value = 913716.48
$ws.Range("N136").Value = -918816.48
